# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.574.41"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.390.47"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'575.63"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'141.35"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'7.67"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").Value = "3.969.92"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "'28.31"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "3.363.67"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "61.608.92"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "'6.12"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'13.57"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "'9.00"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'391.75"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").Value = "'74.93"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'0.0000113"
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("D26").Value = "'0.196"
$ws.Range("E26").Value = "  +8.74%  "
$ws.Range("D27").Value = "'1.01"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Value = "'7.31"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'23.22"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "'168.40"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").Value = "'5.01"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "3.424.43"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").Value = "'0.778"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'4.42"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "'1.15"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "2.459.12"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "'22.61"
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'0.0262"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'2.02"
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("E51").Value = "  -2.00%  "
